# Apply the three changes from the diff:
#  1) Both tables get <w:tblLayout w:type="fixed"/> added to tblPr
#     (AllowAutoFit = $false produces exactly this).
#  2) A new paragraph style "AbstractTitle" ("Abstract Title") is added,
#     based on Normal, next style Abstract.
#  3) The existing "Abstract" style's spacing-before changes from 300 -> 100
#     (twips/20 => 15pt -> 5pt).

$d = $word.ActiveDocument

# --- 1. Tables: force fixed layout (adds w:tblLayout w:type="fixed") ---
foreach ($t in $d.Tables) {
    $t.AllowAutoFit = $false
}

# --- 2. Add the new "AbstractTitle" paragraph style ---
$titleStyle = $d.Styles.Add("AbstractTitle", 1)
$titleStyle.NameLocal = "Abstract Title"
$titleStyle.BaseStyle = "Normal"
$titleStyle.NextParagraphStyle = "Abstract"
$titleStyle.QuickStyle = $true

$titleStyle.ParagraphFormat.KeepWithNext = $true
$titleStyle.ParagraphFormat.KeepTogether = $true
$titleStyle.ParagraphFormat.Alignment = 1
$titleStyle.ParagraphFormat.SpaceAfter = 0
$titleStyle.ParagraphFormat.SpaceBefore = 15

$titleStyle.Font.Size = 10
$titleStyle.Font.SizeBi = 10
$titleStyle.Font.Bold = $true
$titleStyle.Font.Color = 9067060

# --- 3. Abstract style: spacing-before 300 -> 100 (15pt -> 5pt) ---
$abstractStyle = $d.Styles.Item("Abstract")
$abstractStyle.ParagraphFormat.SpaceBefore = 5

Write-Output "done"
